# Update simulation result values in row 2 (results from re-running cases A-F)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -6950.007835040191
$ws.Range("C2").Value = 13573.29778265877
$ws.Range("D2").Value = -6004.444591846527
$ws.Range("E2").Value = -618.8453557719538
$ws.Range("F2").Value = 56.42082461977606
$ws.Range("G2").Value = 45.89295856732645
$ws.Range("H2").Value = 56.68861583280625
$ws.Range("I2").Value = 45.36862496030733
$ws.Range("J2").Value = 56.36079380891869
$ws.Range("K2").Value = 45.21879701254443
$ws.Range("L2").Value = 54.92858262148599
$ws.Range("O2").Value = 45.60957459620175
$ws.Range("P2").Value = 53.57777798835218
$ws.Range("R2").Value = 9.430764574722048
$ws.Range("S2").Value = -17.12936487511296
$ws.Range("T2").Value = 7.698600300390913
$ws.Range("X2").Value = -102.445518799464
$ws.Range("Y2").Value = -157.9190387366542
$ws.Range("Z2").Value = -120.9519319754887
$ws.Range("AE2").Value = -9.430764574722048
$ws.Range("AF2").Value = 7.698600300390913
$ws.Range("AG2").Value = 9.430764574722048
$ws.Range("AH2").Value = -17.12936487511296
$ws.Range("AI2").Value = 7.698600300390913
$ws.Range("AJ2").Value = 9.430764574722048
$ws.Range("AK2").Value = -7.698600300390913
$ws.Range("AL2").Value = 27.73675996859511
$ws.Range("AM2").Value = -18.48355338058273
$ws.Range("AN2").Value = -102.445518799464
$ws.Range("AO2").Value = -157.9190387366542
$ws.Range("AP2").Value = -120.9519319754887
$ws.Range("AQ2").Value = -27.73675996859511
$ws.Range("AR2").Value = 18.48355338058273
$ws.Range("AS2").Value = 56.42082461977606
$ws.Range("AT2").Value = 56.42082461977606
$ws.Range("AU2").Value = 56.68861583280631
$ws.Range("AV2").Value = 56.68861583280631
$ws.Range("AW2").Value = 56.68861583280625
$ws.Range("AX2").Value = 56.36079380891863
$ws.Range("AY2").Value = 56.36079380891869
$ws.Range("AZ2").Value = 45.89295856732645
$ws.Range("BA2").Value = 45.89295856732645
$ws.Range("BB2").Value = 45.36862496030733
$ws.Range("BC2").Value = 45.68985904592068
$ws.Range("BD2").Value = 44.97511407305689
$ws.Range("BE2").Value = 45.21879701254443
$ws.Range("BF2").Value = 45.21879701254443
$ws.Range("BG2").Value = 54.92858262148599
$ws.Range("BJ2").Value = 45.60957459620175
$ws.Range("BK2").Value = 53.57777798835218